# Update the innings-by-innings batting activity table for Prithvi Shaw
# (columns C:runs, D:balls, E:fours, F:sixes) to reflect the latest data
# pulled in from the source (rows were re-ordered/refreshed from the
# Excel activity form).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# r, runs, balls, fours, sixes
$data = @(
    @(2,  9,  6, 2, 0),
    @(3,  10, 11, 2, 0),
    @(4,  0,  2, 0, 0),
    @(5,  0,  1, 0, 0),
    @(6,  0,  2, 0, 0),
    @(7,  4,  3, 1, 0),
    @(8,  7,  11, 1, 0),
    @(9,  2,  5, 0, 0),
    @(10, 66, 41, 4, 4),
    @(11, 64, 43, 9, 1),
    @(12, 5,  9, 1, 0),
    @(13, 42, 23, 5, 2),
    @(14, 19, 10, 2, 1)
)

# Keep these columns formatted/stored as text, matching the source sheet.
$ws.Range("C2:F14").NumberFormat = "@"

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 3).Value = [string]$row[1]
    $ws.Cells.Item($r, 4).Value = [string]$row[2]
    $ws.Cells.Item($r, 5).Value = [string]$row[3]
    $ws.Cells.Item($r, 6).Value = [string]$row[4]
}
